# Apply the commit's change: swap the order of "Recorded By" entries
# from "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"
# for every row where that exact value appears in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $v = $cell.Value()
    if ($v -eq $oldValue) {
        $cell.Value = $newValue
    }
}
